$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the values to be stored as text (not auto-converted to a number
# or a date) by temporarily applying a text number format, then restore
# the original ("Normal") cell style so no formatting change is left behind.
$rangeA = $ws.Range("A2:A15")
$rangeB = $ws.Range("B2:B15")

$rangeA.NumberFormat = "@"
$rangeB.NumberFormat = "@"

$rangeA.Value = "2025"
$rangeB.Value = "June 2025"

$rangeA.Style = "Normal"
$rangeB.Style = "Normal"
